$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Rows 45/46: PaxDollar and Quant swap places (with updated price/volume)
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D45") '103.88'
$ws.Range("E45").Value = '  +0.94%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D46") '0.9996'
$ws.Range("E46").Value = '  -0.12%  '

$ws.Range("D2").Value = '30.302.17'
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").Value = '1.867.50'
$ws.Range("E3").Value = '  -0.66%  '

Set-TextValue $ws.Range("D4") '1.000'

Set-TextValue $ws.Range("D5") '235.43'
$ws.Range("E5").Value = '  +0.21%  '

Set-TextValue $ws.Range("D6") '0.9998'
$ws.Range("E6").Value = '  -0.20%  '

Set-TextValue $ws.Range("D7") '0.4680'
$ws.Range("E7").Value = '  +0.08%  '

Set-TextValue $ws.Range("D8") '0.2843'
$ws.Range("E8").Value = '  +0.67%  '

Set-TextValue $ws.Range("D9") '0.06531'
$ws.Range("E9").Value = '  -0.78%  '

Set-TextValue $ws.Range("D10") '21.44'
$ws.Range("E10").Value = '  +4.30%  '

Set-TextValue $ws.Range("D11") '0.07869'
$ws.Range("E11").Value = '  +1.37%  '

Set-TextValue $ws.Range("D12") '97.82'
$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("D13").Value = '1.865.37'
$ws.Range("E13").Value = '  -0.82%  '

Set-TextValue $ws.Range("D14") '5.103'
$ws.Range("E14").Value = '  +0.79%  '

Set-TextValue $ws.Range("D15") '0.6761'
$ws.Range("E15").Value = '  +0.61%  '

Set-TextValue $ws.Range("D16") '279.49'
$ws.Range("E16").Value = '  -1.42%  '

$ws.Range("D17").Value = '30.294.01'
$ws.Range("E17").Value = '  +0.06%  '

Set-TextValue $ws.Range("D18") '0.9996'
$ws.Range("E18").Value = '  -0.02%  '

Set-TextValue $ws.Range("D19") '5.503'
$ws.Range("E19").Value = '  +2.08%  '

$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("D21").Value = '2.120.26'
$ws.Range("E21").Value = '  -0.24%  '

Set-TextValue $ws.Range("D22") '0.000007293'
$ws.Range("E22").Value = '  +0.69%  '

Set-TextValue $ws.Range("D23") '0.9994'
$ws.Range("E23").Value = '  -0.23%  '

Set-TextValue $ws.Range("D24") '6.168'
$ws.Range("E24").Value = '  +0.11%  '

Set-TextValue $ws.Range("D25") '9.189'
$ws.Range("E25").Value = '  -1.68%  '

$ws.Range("E26").Value = '  -1.59%  '

$ws.Range("E27").Value = '  +0.00%  '

Set-TextValue $ws.Range("D28") '1.930'
$ws.Range("E28").Value = '  -2.59%  '

Set-TextValue $ws.Range("D29") '1.375'
$ws.Range("E29").Value = '  -0.47%  '

Set-TextValue $ws.Range("D30") '0.09652'
$ws.Range("E30").Value = '  -0.12%  '

Set-TextValue $ws.Range("D31") '4.379'
$ws.Range("E31").Value = '  +0.42%  '

$ws.Range("E32").Value = '  +0.30%  '

Set-TextValue $ws.Range("D33") '4.101'
$ws.Range("E33").Value = '  -0.03%  '

Set-TextValue $ws.Range("D34") '0.04712'
$ws.Range("E34").Value = '  +1.14%  '

Set-TextValue $ws.Range("D35") '1.129'
$ws.Range("E35").Value = '  +3.28%  '

Set-TextValue $ws.Range("D36") '0.7066'
$ws.Range("E36").Value = '  +0.41%  '

Set-TextValue $ws.Range("D37") '2.719'
$ws.Range("E37").Value = '  +0.18%  '

Set-TextValue $ws.Range("D38") '0.01856'
$ws.Range("E38").Value = '  -0.62%  '

Set-TextValue $ws.Range("D39") '6.278'
$ws.Range("E39").Value = '  -4.50%  '

Set-TextValue $ws.Range("D40") '2.527'
$ws.Range("E40").Value = '  +0.12%  '

Set-TextValue $ws.Range("D41") '73.68'
$ws.Range("E41").Value = '  +2.58%  '

Set-TextValue $ws.Range("D42") '1.948'
$ws.Range("E42").Value = '  -0.35%  '

Set-TextValue $ws.Range("D43") '0.8480'
$ws.Range("E43").Value = '  -1.68%  '

Set-TextValue $ws.Range("D44") '0.4179'
$ws.Range("E44").Value = '  +0.16%  '

Set-TextValue $ws.Range("D47") '7.177'
$ws.Range("E47").Value = '  -0.93%  '

Set-TextValue $ws.Range("D48") '9.182'
$ws.Range("E48").Value = '  +0.36%  '

Set-TextValue $ws.Range("D49") '935.72'
$ws.Range("E49").Value = '  -4.73%  '

Set-TextValue $ws.Range("D50") '34.12'
$ws.Range("E50").Value = '  +0.93%  '

Set-TextValue $ws.Range("D51") '0.1125'
$ws.Range("E51").Value = '  -1.60%  '

